# Updates the cryptos price list (columns B-E) to the refreshed
# GitHub Actions snapshot described by the commit diff.
#
# Set-CellText writes a value as literal text even when it "looks"
# numeric (e.g. "213.90", "1.000"), matching the original
# <is><t>...</t></is> inline-string cells -- a bare $range.Value =
# assignment lets Excel auto-coerce those into real numbers and
# silently rewrite "1.000" -> 1, "19.60" -> 19.6, etc.
function Set-CellText($range, $text) {
    $range.Value = "'" + $text
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.795.08'
$ws.Range("E2").Value = '  -1.28%  '
$ws.Range("D3").Value = '1.632.02'
$ws.Range("E3").Value = '  -1.23%  '
$ws.Range("E4").Value = '  -0.13%  '
Set-CellText $ws.Range("D5") '213.90'
$ws.Range("E5").Value = '  -0.83%  '
Set-CellText $ws.Range("D6") '0.5015'
$ws.Range("E6").Value = '  -1.69%  '
Set-CellText $ws.Range("D7") '1.000'
$ws.Range("E7").Value = '  -0.48%  '
Set-CellText $ws.Range("D8") '0.2554'
$ws.Range("E8").Value = '  -0.95%  '
Set-CellText $ws.Range("D9") '0.06358'
$ws.Range("E9").Value = '  -0.79%  '
Set-CellText $ws.Range("D10") '19.60'
$ws.Range("E10").Value = '  -1.65%  '
Set-CellText $ws.Range("D11") '0.07707'
$ws.Range("E11").Value = '  -1.06%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.652.42'
$ws.Range("E12").Value = '  -0.43%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-CellText $ws.Range("D13") '4.251'
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("D14").Value = '1.858.94'
$ws.Range("E14").Value = '  -1.09%  '
Set-CellText $ws.Range("D15") '0.5410'
$ws.Range("E15").Value = '  -1.90%  '
$val = '0.0{0}7865' -f [char]0x2085
Set-CellText $ws.Range("D16") $val
$ws.Range("E16").Value = '  -1.78%  '
Set-CellText $ws.Range("D17") '63.50'
$ws.Range("E17").Value = '  -0.81%  '
$ws.Range("D18").Value = '25.809.89'
$ws.Range("E18").Value = '  -1.20%  '
Set-CellText $ws.Range("D19") '1.001'
$ws.Range("E19").Value = '  -0.50%  '
Set-CellText $ws.Range("D20") '200.10'
$ws.Range("E20").Value = '  -4.66%  '
Set-CellText $ws.Range("D21") '4.320'
$ws.Range("E21").Value = '  -2.14%  '
Set-CellText $ws.Range("D22") '9.851'
$ws.Range("E22").Value = '  -2.16%  '
Set-CellText $ws.Range("D23") '5.927'
$ws.Range("E23").Value = '  -1.86%  '
Set-CellText $ws.Range("D24") '1.001'
$ws.Range("E24").Value = '  -0.45%  '
Set-CellText $ws.Range("D25") '1.918'
$ws.Range("E25").Value = '  +10.24%  '
Set-CellText $ws.Range("D26") '140.70'
$ws.Range("E26").Value = '  -2.06%  '
Set-CellText $ws.Range("D27") '0.1131'
$ws.Range("E27").Value = '  -4.05%  '
Set-CellText $ws.Range("D28") '15.60'
$ws.Range("E28").Value = '  -1.44%  '
Set-CellText $ws.Range("D29") '6.664'
$ws.Range("E29").Value = '  -4.37%  '
Set-CellText $ws.Range("D30") '1.235'
$ws.Range("E30").Value = '  -0.49%  '
Set-CellText $ws.Range("D31") '0.04963'
$ws.Range("E31").Value = '  -2.57%  '
Set-CellText $ws.Range("D32") '3.258'
$ws.Range("E32").Value = '  -2.51%  '
Set-CellText $ws.Range("D33") '3.174'
$ws.Range("E33").Value = '  -1.37%  '
Set-CellText $ws.Range("D34") '1.530'
$ws.Range("E34").Value = '  -2.30%  '
Set-CellText $ws.Range("D35") '2.362'
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("B36").Value = 'MXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-CellText $ws.Range("D36") '2.627'
$ws.Range("E36").Value = '  -4.41%  '
$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '1.163.41'
$ws.Range("E37").Value = '  -0.23%  '
Set-CellText $ws.Range("D38") '0.8837'
$ws.Range("E38").Value = '  -4.50%  '
Set-CellText $ws.Range("D39") '0.5545'
$ws.Range("E39").Value = '  -2.31%  '
Set-CellText $ws.Range("D40") '0.01554'
$ws.Range("E40").Value = '  -2.20%  '
Set-CellText $ws.Range("D41") '0.9989'
$ws.Range("E41").Value = '  -0.61%  '
Set-CellText $ws.Range("D42") '5.639'
$ws.Range("E42").Value = '  -0.30%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-CellText $ws.Range("D43") '99.14'
$ws.Range("E43").Value = '  -1.26%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-CellText $ws.Range("D44") '0.8007'
$ws.Range("E44").Value = '  -3.57%  '
$ws.Range("D45").Value = '1.771.72'
$ws.Range("E45").Value = '  -1.03%  '
$ws.Range("E46").Value = '  +0.04%  '
Set-CellText $ws.Range("D47") '0.4519'
$ws.Range("E47").Value = '  -0.72%  '
Set-CellText $ws.Range("D48") '1.002'
$ws.Range("E48").Value = '  -0.43%  '
Set-CellText $ws.Range("D49") '54.42'
$ws.Range("E49").Value = '  -2.19%  '
Set-CellText $ws.Range("D50") '0.05061'
$ws.Range("E50").Value = '  +0.05%  '
Set-CellText $ws.Range("D51") '1.002'
$ws.Range("E51").Value = '  -0.45%  '
